$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.547.50'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '1.663.53'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.80'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4790'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2612'
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07086'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').Value = '1.665.92'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.72'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.5901'
$ws.Range('E13').Value = '  -4.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.366'
$ws.Range('E14').Value = '  -4.35%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '74.31'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '25.538.21'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000006749'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '1.879.96'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.423'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.650'
$ws.Range('E23').Value = '  +2.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.293'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '134.46'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.03'
$ws.Range('E26').Value = '  +2.04%  '
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('E28').Value = '  +2.62%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.683'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.960'
$ws.Range('E30').Value = '  +4.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.650'
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('E32').Value = '  -3.85%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9996'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  -5.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.616'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6111'
$ws.Range('E36').Value = '  +5.63%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9473'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.8477'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.01493'
$ws.Range('E41').Value = '  -2.90%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.865'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '97.95'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3755'
$ws.Range('E44').Value = '  +1.49%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.201'
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.41'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('E51').Value = '  +0.36%  '
